$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.341.09'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = '2.938.35'
$ws.Range('E3').Value = '  +3.87%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '200.08'
$ws.Range('E5').Value = '  +4.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '596.99'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.552'
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.196'
$ws.Range('D10').Value = '2.940.06'
$ws.Range('E10').Value = '  +3.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.447'
$ws.Range('E11').Value = '  +16.75%  '
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').Value = '3.480.15'
$ws.Range('E14').Value = '  +3.80%  '
$ws.Range('D15').Value = '76.265.10'
$ws.Range('E15').Value = '  +1.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '28.04'
$ws.Range('E16').Value = '  +3.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000189'
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('D18').Value = '2.910.36'
$ws.Range('E18').Value = '  +3.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.33'
$ws.Range('E19').Value = '  +8.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.72'
$ws.Range('E20').Value = '  -3.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '373.26'
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.30'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('E23').Value = '  +5.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.65'
$ws.Range('E24').Value = '  +2.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').Value = '3.057.85'
$ws.Range('E26').Value = '  +3.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.28'
$ws.Range('E27').Value = '  +2.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.69'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000108'
$ws.Range('E29').Value = '  +3.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('E31').Value = '  -2.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.90'
$ws.Range('E32').Value = '  +2.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '497.75'
$ws.Range('E33').Value = '  -3.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.83'
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '165.27'
$ws.Range('E36').Value = '  +0.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.26'
$ws.Range('E37').Value = '  +1.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.395'
$ws.Range('E38').Value = '  +15.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.109'
$ws.Range('E39').Value = '  +24.76%  '
$ws.Range('E40').Value = '  +1.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.111'
$ws.Range('E41').Value = '  -5.70%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '178.63'
$ws.Range('E43').Value = '  -2.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.94'
$ws.Range('E44').Value = '  -1.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.65'
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.21'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('E47').Value = '  -1.33%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.89'
$ws.Range('E48').Value = '  +3.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.584'
$ws.Range('E49').Value = '  +2.29%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.31'
$ws.Range('E50').Value = '  -2.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.76'
$ws.Range('E51').Value = '  +7.88%  '
